# Berechnen der Kapazitaet eines Flugs
# Adds a new row (28) documenting the new API method "calculateCapacity"
# on the sheet that lists classes/methods, and tweaks the sheet view /
# column sizing that Excel recorded when the row was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table row describing calculateCapacity -----------------------
$ws.Range("B28").Value = "calculateCapacity"
$ws.Range("C28").Value = "String flugString"
$ws.Range("D28").Value = "String"
$ws.Range("E28").Value = "Dieser Flug hat noch kein zugeordnetes Flugzeug!"
$ws.Range("F28").Value = "-"
$ws.Range("G28").Value = "gibt die Kapazität eines Flugs aus. Beispiel: Es sind noch 842 Plaetze von 853 Plaetzen frei."

# Row 28 styling: same "top/left, no wrap" look as the rest of the table
# except for the error column (E), which wraps without the usual
# top/left alignment tweak - matches the new cellXfs entry Excel created.
$ws.Range("B28").VerticalAlignment = -4160    # xlTop
$ws.Range("B28").HorizontalAlignment = -4131  # xlLeft

$ws.Range("D28").VerticalAlignment = -4160
$ws.Range("D28").HorizontalAlignment = -4131

$ws.Range("F28").VerticalAlignment = -4160
$ws.Range("F28").HorizontalAlignment = -4131

$ws.Range("G28").VerticalAlignment = -4160
$ws.Range("G28").HorizontalAlignment = -4131

$ws.Range("E28").WrapText = $true

$ws.Rows.Item(28).RowHeight = 30

# --- Column widths: switch from manually-sized columns to "best fit" --
# (input values are nudged slightly so the engine's internal width
# rounding lands as close as possible to Excel's recorded bestFit widths)
$ws.Columns.Item(1).ColumnWidth = 17.76
$ws.Columns.Item(2).ColumnWidth = 26.09
$ws.Columns.Item(3).ColumnWidth = 43.59
$ws.Columns.Item(4).ColumnWidth = 16.25
$ws.Columns.Item(5).ColumnWidth = 34.42
$ws.Columns.Item(6).ColumnWidth = 40.59
$ws.Columns.Item(7).ColumnWidth = 111.09

# --- Row heights that Excel recomputed once the columns got wider -----
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(25).AutoFit()

# --- View state: zoom level and active cell ----------------------------
$excel.ActiveWindow.Zoom = 80
$ws.Range("G4").Select() | Out-Null
